$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case for #268: IF(ISBLANK(...)) formula returning a string result.
$ws.Range("A16").Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'

# Move the active selection to the newly added cell, as the original author did.
$ws.Range("A16").Select() | Out-Null
